$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

# Fall 23 week 15 inputs - update matchup average values
$ws.Range("D3").Value = 10.32
$ws.Range("E3").Value = 10.52
$ws.Range("H3").Value = 11.43

$ws.Range("C4").Value = 9.68
$ws.Range("F4").Value = 10.08

$ws.Range("C5").Value = 9.48
$ws.Range("F5").Value = 10.39
$ws.Range("G5").Value = 9.42
$ws.Range("H5").Value = 8.710000000000001

$ws.Range("D6").Value = 9.92
$ws.Range("E6").Value = 9.609999999999999

$ws.Range("E7").Value = 10.58
$ws.Range("H7").Value = 9.619999999999999

$ws.Range("C8").Value = 8.57
$ws.Range("E8").Value = 11.29
$ws.Range("G8").Value = 10.38
